$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Date" value (row 8, column B) to reflect the new generation timestamp.
$ws.Range("B8").Value = "2022-01-21T07:49:24+01:00"

# The profile previously listed two "Contact" rows (rows 10-11), each holding
# "No display for ContactDetail". The IG now lists all authors as contacts,
# which doubles the contact rows to four (rows 10-13). Shift every row below
# the existing contact rows down by two to make room, then populate the two
# new contact rows.

for ($r = 21; $r -ge 12; $r--) {
    $dest = $r + 2
    $ws.Range("A$dest").ClearContents()
    $ws.Range("A$r").Copy($ws.Range("A$dest"))
    $ws.Range("B$dest").ClearContents()
    $ws.Range("B$r").Copy($ws.Range("B$dest"))
}

$ws.Range("A12").ClearContents()
$ws.Range("A10").Copy($ws.Range("A12"))
$ws.Range("B12").ClearContents()
$ws.Range("B10").Copy($ws.Range("B12"))

$ws.Range("A13").ClearContents()
$ws.Range("A10").Copy($ws.Range("A13"))
$ws.Range("B13").ClearContents()
$ws.Range("B10").Copy($ws.Range("B13"))
